$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (prices) that must remain text,
# matching the source file's inlineStr cells. We briefly force a text number
# format before assigning so Excel doesn't auto-convert the string to a number,
# then restore the default 'Normal' style so no stray formatting is left behind.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "93.737.81"
$ws.Range("E2").Value = "  +1.33%  "

Set-TextValue $ws.Range("D3") "3.086.53"
$ws.Range("E3").Value = "  -0.82%  "

Set-TextValue $ws.Range("D5") "233.59"
$ws.Range("E5").Value = "  -3.32%  "

Set-TextValue $ws.Range("D6") "609.12"
$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("E8").Value = "  -5.49%  "

$ws.Range("E9").Value = "  +0.02%  "

Set-TextValue $ws.Range("D10") "0.817"
$ws.Range("E10").Value = "  +11.62%  "

Set-TextValue $ws.Range("D11") "3.089.72"
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("E12").Value = "  -3.57%  "

Set-TextValue $ws.Range("D13") "93.867.76"
$ws.Range("E13").Value = "  +1.73%  "

$ws.Range("E14").Value = "  -6.66%  "

Set-TextValue $ws.Range("D15") "34.01"
$ws.Range("E15").Value = "  -1.16%  "

Set-TextValue $ws.Range("D16") "3.666.81"
$ws.Range("E16").Value = "  -0.60%  "

Set-TextValue $ws.Range("D17") "5.24"
$ws.Range("E17").Value = "  -4.94%  "

Set-TextValue $ws.Range("D18") "3.128.32"
$ws.Range("E18").Value = "  +1.16%  "

Set-TextValue $ws.Range("D19") "3.64"
$ws.Range("E19").Value = "  -0.90%  "

Set-TextValue $ws.Range("D20") "14.63"
$ws.Range("E20").Value = "  -0.87%  "

Set-TextValue $ws.Range("D21") "5.76"
$ws.Range("E21").Value = "  -0.70%  "

Set-TextValue $ws.Range("D22") "439.59"
$ws.Range("E22").Value = "  -1.87%  "

Set-TextValue $ws.Range("D23") "8.83"
$ws.Range("E23").Value = "  -6.25%  "

Set-TextValue $ws.Range("D24") "0.0000191"
$ws.Range("E24").Value = "  -5.78%  "

Set-TextValue $ws.Range("D25") "8.27"
$ws.Range("E25").Value = "  +4.91%  "

Set-TextValue $ws.Range("D26") "5.51"
$ws.Range("E26").Value = "  -4.43%  "

Set-TextValue $ws.Range("D27") "84.65"
$ws.Range("E27").Value = "  -2.98%  "

Set-TextValue $ws.Range("D28") "11.96"
$ws.Range("E28").Value = "  +1.65%  "

Set-TextValue $ws.Range("D29") "3.269.98"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("E30").Value = "  -0.02%  "

Set-TextValue $ws.Range("D31") "0.252"
$ws.Range("E31").Value = "  +8.63%  "

Set-TextValue $ws.Range("D32") "0.178"
$ws.Range("E32").Value = "  +5.70%  "

Set-TextValue $ws.Range("D33") "0.122"
$ws.Range("E33").Value = "  -11.05%  "

Set-TextValue $ws.Range("D34") "9.22"
$ws.Range("E34").Value = "  -0.81%  "

$ws.Range("E35").Value = "  +0.40%  "

Set-TextValue $ws.Range("D36") "7.73"
$ws.Range("E36").Value = "  -4.20%  "

$ws.Range("E37").Value = "  -4.92%  "

Set-TextValue $ws.Range("D38") "25.53"
$ws.Range("E38").Value = "  -2.46%  "

$ws.Range("E39").Value = "  -1.94%  "

Set-TextValue $ws.Range("D40") "0.444"
$ws.Range("E40").Value = "  +0.36%  "

$ws.Range("E41").Value = "  +3.69%  "

Set-TextValue $ws.Range("D42") "1.28"
$ws.Range("E42").Value = "  -2.69%  "

$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Range("D43") "3.68"
$ws.Range("E43").Value = "  -13.64%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D44") "464.12"
$ws.Range("E44").Value = "  -3.57%  "

$ws.Range("E45").Value = "  -0.01%  "

Set-TextValue $ws.Range("D46") "3.10"
$ws.Range("E46").Value = "  -11.30%  "

Set-TextValue $ws.Range("D47") "161.42"

$ws.Range("E48").Value = "  -1.96%  "

$ws.Range("E49").Value = "  -3.61%  "

Set-TextValue $ws.Range("D50") "43.70"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("E51").Value = "  -0.04%  "
